# Auto-generated Excel COM-interop script
# Updates FFXIV Leve-crafting market-price columns (H-N) across all class sheets
# to reflect the latest scraped market data, per the scheduled runner job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 250.26923
$ws.Range("I33").Value = 194.42857
$ws.Range("K33").Value = 194.42857
$ws.Range("M33").Value = 34.57142999999999

$ws.Range("H70").Value = 1705.5
$ws.Range("I70").Value = 1436.3334
$ws.Range("K70").Value = 4309.0002
$ws.Range("M70").Value = -4039.0002

$ws.Range("H73").Value = 1705.5
$ws.Range("I73").Value = 1436.3334
$ws.Range("K73").Value = 4309.0002
$ws.Range("M73").Value = -3373.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3353.25
$ws.Range("I32").Value = 2940.463
$ws.Range("K32").Value = 2940.463
$ws.Range("M32").Value = -2653.463

$ws.Range("H45").Value = 1240.92
$ws.Range("I45").Value = 1163.5883
$ws.Range("J45").Value = 1405.25
$ws.Range("K45").Value = 1163.5883
$ws.Range("L45").Value = 1405.25
$ws.Range("M45").Value = -786.5882999999999
$ws.Range("N45").Value = -2159.25

$ws.Range("H61").Value = 3614.8276
$ws.Range("I61").Value = 3442.5652
$ws.Range("K61").Value = 3442.5652
$ws.Range("M61").Value = -3230.5652

$ws.Range("H122").Value = 6459835.5
$ws.Range("I122").Value = 8336922
$ws.Range("K122").Value = 25010766
$ws.Range("M122").Value = -25008316

$ws.Range("H132").Value = 6623.989
$ws.Range("I132").Value = 4420.3115
$ws.Range("K132").Value = 13260.9345
$ws.Range("M132").Value = -10730.9345

$ws.Range("H136").Value = 3614.8276
$ws.Range("I136").Value = 3442.5652
$ws.Range("K136").Value = 10327.6956
$ws.Range("M136").Value = -7777.695599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 38463560
$ws.Range("I134").Value = 45456652
$ws.Range("J134").Value = 1555.5
$ws.Range("K134").Value = 136369956
$ws.Range("L134").Value = 4666.5
$ws.Range("M134").Value = -136367421
$ws.Range("N134").Value = -9736.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1855.8928
$ws.Range("I31").Value = 1711.5
$ws.Range("K31").Value = 1711.5
$ws.Range("M31").Value = -1416.5

$ws.Range("H34").Value = 1855.8928
$ws.Range("I34").Value = 1711.5
$ws.Range("K34").Value = 1711.5
$ws.Range("M34").Value = -1509.5

$ws.Range("H39").Value = 14666.333
$ws.Range("I39").Value = 14666.333
$ws.Range("K39").Value = 14666.333
$ws.Range("M39").Value = -14275.333

$ws.Range("H49").Value = 14666.333
$ws.Range("I49").Value = 14666.333
$ws.Range("K49").Value = 14666.333
$ws.Range("M49").Value = -14484.333

$ws.Range("H99").Value = 2824.375
$ws.Range("I99").Value = 2332.3333
$ws.Range("K99").Value = 2332.3333
$ws.Range("M99").Value = -834.3332999999998

$ws.Range("H122").Value = 2703.4119
$ws.Range("I122").Value = 1703.9
$ws.Range("K122").Value = 5111.700000000001
$ws.Range("M122").Value = -2661.700000000001

$ws.Range("H126").Value = 2824.375
$ws.Range("I126").Value = 2332.3333
$ws.Range("K126").Value = 6996.999899999999
$ws.Range("M126").Value = -4526.999899999999

$ws.Range("H132").Value = 1146.1333
$ws.Range("I132").Value = 1156.5714
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3469.7142
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -939.7142000000003
$ws.Range("N132").Value = -8060

$ws.Range("H134").Value = 1347.96
$ws.Range("I134").Value = 1281.8182
$ws.Range("K134").Value = 3845.4546
$ws.Range("M134").Value = -1310.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 137.3077
$ws.Range("J12").Value = 134.88889
$ws.Range("L12").Value = 404.66667
$ws.Range("N12").Value = -750.6666700000001

$ws.Range("H39").Value = 750
$ws.Range("I39").Value = 750
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1956
$ws.Range("N39").Value = ""

$ws.Range("H69").Value = 10010.25
$ws.Range("I69").Value = 10010.25
$ws.Range("K69").Value = 30030.75
$ws.Range("M69").Value = -29219.75

$ws.Range("H72").Value = 10010.25
$ws.Range("I72").Value = 10010.25
$ws.Range("K72").Value = 90092.25
$ws.Range("M72").Value = -86036.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2535.862
$ws.Range("I132").Value = 1667.75
$ws.Range("J132").Value = 3604.3076
$ws.Range("K132").Value = 5003.25
$ws.Range("L132").Value = 10812.9228
$ws.Range("M132").Value = -2473.25
$ws.Range("N132").Value = -15872.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 254249.75
$ws.Range("I22").Value = 1000000
$ws.Range("J22").Value = 5666.3335
$ws.Range("K22").Value = 1000000
$ws.Range("L22").Value = 5666.3335
$ws.Range("M22").Value = -999705
$ws.Range("N22").Value = -6256.3335

$ws.Range("H27").Value = 254249.75
$ws.Range("I27").Value = 1000000
$ws.Range("J27").Value = 5666.3335
$ws.Range("K27").Value = 1000000
$ws.Range("L27").Value = 5666.3335
$ws.Range("M27").Value = -999893
$ws.Range("N27").Value = -5880.3335

$ws.Range("H51").Value = 40084
$ws.Range("J51").Value = 40084
$ws.Range("L51").Value = 40084
$ws.Range("N51").Value = -41040

$ws.Range("H55").Value = 973.5625
$ws.Range("I55").Value = 416
$ws.Range("J55").Value = 1788.4615
$ws.Range("K55").Value = 416
$ws.Range("L55").Value = 1788.4615
$ws.Range("M55").Value = -243
$ws.Range("N55").Value = -2134.4615

$ws.Range("H61").Value = 943.4167
$ws.Range("I61").Value = 847.36365
$ws.Range("K61").Value = 847.36365
$ws.Range("M61").Value = -645.36365

$ws.Range("H68").Value = 7427
$ws.Range("J68").Value = 11498
$ws.Range("L68").Value = 11498
$ws.Range("N68").Value = -12996

$ws.Range("H71").Value = 7427
$ws.Range("J71").Value = 11498
$ws.Range("L71").Value = 57490
$ws.Range("N71").Value = -64978

$ws.Range("H113").Value = 943.4167
$ws.Range("I113").Value = 847.36365
$ws.Range("K113").Value = 847.36365
$ws.Range("M113").Value = 1322.63635

$ws.Range("H132").Value = 3840.0435
$ws.Range("I132").Value = 3488.7058
$ws.Range("K132").Value = 10466.1174
$ws.Range("M132").Value = -7936.117400000001

$ws.Range("H136").Value = 2866.7144
$ws.Range("I136").Value = 2259.1628
$ws.Range("K136").Value = 6777.4884
$ws.Range("M136").Value = -4227.4884

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3358.875
$ws.Range("I122").Value = 2876.8635
$ws.Range("K122").Value = 8630.5905
$ws.Range("M122").Value = -6180.5905

$ws.Range("H126").Value = 1919.6
$ws.Range("I126").Value = 1744
$ws.Range("K126").Value = 5232
$ws.Range("M126").Value = -2762
